$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: Capture formatting for the "style-only" helper cells from
# stable reference cells elsewhere on the sheet (style index 2 and 3)
# before we start mutating rows 52-54.
# ------------------------------------------------------------------
# (We copy directly from L48 -> style "2", and H8 -> style "3" later.)

# ------------------------------------------------------------------
# Step 2: Overwrite row 52 (currently "Kf") with what used to be row
# 53's data ("Ks"), keeping the same per-cell styles row 53 had.
# ------------------------------------------------------------------
$ws.Range("B53").Copy()
$ws.Range("B52").PasteSpecial(-4122)
$ws.Range("A52").Value = "Ks"
$ws.Range("B52").Value = 360
$ws.Range("C52").Value = "Infiltration rate of submerged zone m/hr"
$ws.Range("D52").Value = "Free draining"

# Clear any stray leftover cells from the old row 52 content (I52) that
# are not part of the new row 52 content.
$ws.Range("I52").Clear()

# G52 needs the blank styled placeholder (style "3") that used to sit
# at G53.
$ws.Range("H8").Copy()
$ws.Range("G52").PasteSpecial(-4122)

# L52 keeps the blank styled placeholder (style "2"); already present,
# but re-assert it from a stable source for safety.
$ws.Range("L48").Copy()
$ws.Range("L52").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 3: Remove rows 53 and 54 entirely (they become blank rows with
# no cells at all, so they disappear from the sheet XML, and every
# row below keeps its original row number).
# ------------------------------------------------------------------
$ws.Range("A53:S54").Clear()

# ------------------------------------------------------------------
# Step 4: Update row 77 ("native_depth") value + style.
# ------------------------------------------------------------------
$ws.Range("B42").Copy()
$ws.Range("B77").PasteSpecial(-4122)
$ws.Range("B77").Value = 0.19768124000000001

# ------------------------------------------------------------------
# Step 5: Append new row 78 with the old "Kf" row content (now using a
# newly re-calibrated value) below the last row.
# ------------------------------------------------------------------
$ws.Range("A78").Value = "Kf"
$ws.Range("B8").Copy()
$ws.Range("B78").PasteSpecial(-4122)
$ws.Range("B78").Value = 0.11874427999999999
$ws.Range("C78").Value = " Infiltration rate of filter zone m/hr. Calibrated 0.40216883, assumed 1.2 to start based on water-balance"
$ws.Range("D78").Value = "Calibrated value from flows "

$ws.Range("H8").Copy()
$ws.Range("I78").PasteSpecial(-4122)

$ws.Range("L48").Copy()
$ws.Range("L78").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 6: Append new row 79 with the old "Kn" row content (now using a
# newly re-calibrated value) below row 78.
# ------------------------------------------------------------------
$ws.Range("A79").Value = "Kn"
$ws.Range("B79").Value = 0.19605502
$ws.Range("C79").Value = "Saturated infiltration rate of native soil m/hr . "
$ws.Range("D79").Value = "Typical value for S. Ontario silty soil https://wiki.sustainabletechnologies.ca/index.php?title=Low_permeability_soils&mobileaction=toggle_view_desktop 3.3e-3, calibrated from flows"

$ws.Range("L48").Copy()
$ws.Range("L79").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 7: Update the sheet view (selected cell).
# ------------------------------------------------------------------
$ws.Range("C78").Select()
